$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53.. down by one.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new data record.
$ws.Cells.Item(53, 1).Value = 6
$ws.Cells.Item(53, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(53, 3).Value = "Metropolitana"
$ws.Cells.Item(53, 4).Value = 44536
$ws.Cells.Item(53, 5).Value = 13
$ws.Cells.Item(53, 6).Value = 100112001
$ws.Cells.Item(53, 7).Value = "Berenjena"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 210
$ws.Cells.Item(53, 11).Value = 9500
$ws.Cells.Item(53, 12).Value = 10000
$ws.Cells.Item(53, 13).Value = 9786
$ws.Cells.Item(53, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(53, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 16).Value = 196
$ws.Cells.Item(53, 17).Value = 50
$ws.Cells.Item(53, 18).Value = "Hortaliza"
